# إضافة حدث جديد في Card24
# The sheet already had a service-event row (row 20) whose "empty" columns
# were stored as blank cells. This edit:
#   1. Duplicates that existing row down to a new row 21 (preserving the
#      original blank cells as-is, and keeping text types intact).
#   2. Normalizes the previously blank cells of row 20 to the placeholder
#      text "nan", matching the pattern used by every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Step 1: copy the whole existing event row (20) down to the new row (21) ---
# Use Copy + PasteSpecial(values) so the new row's cells keep their original
# data types (e.g. "23" stays text) instead of Excel's normal "typed value"
# auto-detection that a plain .Value assignment of a numeric-looking string
# would trigger.
$ws.Range("A20:P20").Copy()
$ws.Range("A21").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Step 2: fill the blank cells in row 20 with the "nan" placeholder ---
$ws.Range("B20").Value = "nan"
$ws.Range("C20").Value = "nan"
$ws.Range("D20").Value = "nan"
$ws.Range("E20").Value = "nan"
$ws.Range("F20").Value = "nan"
$ws.Range("G20").Value = "nan"
$ws.Range("H20").Value = "nan"
$ws.Range("I20").Value = "nan"
$ws.Range("J20").Value = "nan"
$ws.Range("K20").Value = "nan"
$ws.Range("M20").Value = "nan"
$ws.Range("P20").Value = "nan"
